# A new weekly price record was inserted into the "Poroto verde" price
# table at row 121 (pushing the existing rows 121-215 down to 122-216).
# Replicate that with a real row insert followed by populating the new
# row's cells with the recorded values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 121; everything at/after
# row 121 (including formatting) shifts down by one row.
$ws.Rows.Item(121).Insert()

# Fill in the new row 121 with the new market record.
$ws.Range("A121").Value = 8
$ws.Range("B121").Value = "Terminal La Palmera de La Serena"
$ws.Range("C121").Value = "Coquimbo"
$ws.Range("D121").Value = 44673
$ws.Range("E121").Value = 4
$ws.Range("F121").Value = 100112031
$ws.Range("G121").Value = "Poroto verde"
$ws.Range("H121").Value = "Magnum"
$ws.Range("I121").Value = "Primera"
$ws.Range("J121").Value = 500
$ws.Range("K121").Value = 19000
$ws.Range("L121").Value = 20000
$ws.Range("M121").Value = 19500
$ws.Range("N121").Value = "`$/malla 25 kilos"
$ws.Range("O121").Value = "Provincia de Limarí"
$ws.Range("P121").Value = 780
$ws.Range("Q121").Value = 25
$ws.Range("R121").Value = "Hortaliza"
